# BOM update: add sockets for IMU, Sensors and ODrives.
# - Insert a new "Product" column (E) between Name/Type (D) and price (old E).
# - Rename the "Name" header to "Type".
# - Add a "Digikey" link column after Mouser.
# - Add "ordered"/"stock" tracking columns (J/K).
# - Add new BOM rows for R3 (shunt) and C1 (LDO regulator).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new column E ("Product") ---------------------------------
# Shifts old E (price, width 6) -> F and old F (Mouser url, width 114.7) -> G
$ws.Columns.Item(5).Insert()
$ws.Columns.Item(5).ColumnWidth = 17.59

# --- Header row ---------------------------------------------------------
$ws.Range("D1").Value = "Type"
$ws.Range("E1").Value = "Product"
$ws.Range("H1").Value = "Digikey"
$ws.Range("J1").Value = "ordered"
$ws.Range("K1").Value = "stock"

# --- Existing rows: add Digikey links + ordered counts -----------------
$ws.Range("H2").Value = "https://www.digikey.de/de/products/detail/stmicroelectronics/STL100N10F7/3993110?utm_medium=email&utm_source=oce&utm_campaign=3310_OCE22RT&utm_content=productdetail_DE&utm_cid=2355993&so=75928733&mkt_tok=MDI4LVNYSy01MDcAAAGEemh6HtzLz_P-FEsSsUaBfWxGdHeMKJiNMf369mj9_Qi4I8NoYL4rC2tRRaQ6gPFC7QiGalt27_JUSQXtyGJ6SEnBeNibavWHfdGxt7tI"
$ws.Range("J2").Value = 4

$ws.Range("H3").Value = "https://www.digikey.de/de/products/detail/analog-devices-inc/LTC7001IMSE-PBF/7363804?utm_medium=email&utm_source=oce&utm_campaign=3310_OCE22RT&utm_content=productdetail_DE&utm_cid=2355993&so=75928733&mkt_tok=MDI4LVNYSy01MDcAAAGEemh6Hpx5pIinyDJpdIZawselVf67jNKToa_lz8TVHDCvkXFgNdSQy3jIBmqhnW0PfZFyGL8vHVQNrkMXCLy7Lbft8ti-Gt_-2fDYM6qN"
$ws.Range("J3").Value = 4

$ws.Range("H4").Value = "https://www.digikey.de/de/products/detail/analog-devices-inc-maxim-integrated/MAX3222ECWN-T/1514556?utm_medium=email&utm_source=oce&utm_campaign=3310_OCE22RT&utm_content=productdetail_DE&utm_cid=2355993&so=75928733&mkt_tok=MDI4LVNYSy01MDcAAAGEemh6Hn6FuIg2MPNw_UjWh-MqLqpeDh-4ApnQN1rfi1xcU3gAKpXXwtZ6H_xmgn_GxLTfSQwX3hPh8praGgpQrDjMwOVEJ7bu-RNTSunv"
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1

# --- New row 5: R3 shunt resistor ---------------------------------------
$ws.Range("A5").Value = "Central board"
$ws.Range("B5").Value = "R3"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "Shunt 0.0004Ohm"

# --- New row 6: C1 LDO regulator -----------------------------------------
$ws.Range("A6").Value = "Central board"
$ws.Range("B6").Value = "C1"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = "LDO Regulator 3.3.V"
$ws.Range("E6").Value = "AP2128K-3.3TRG1"
$ws.Range("H6").Value = "https://www.digikey.de/de/products/detail/diodes-incorporated/AP2128K-3-3TRG1/4470794?utm_medium=email&utm_source=oce&utm_campaign=3310_OCE22RT&utm_content=productdetail_DE&utm_cid=2355993&so=75928733&mkt_tok=MDI4LVNYSy01MDcAAAGEemh6HiH8Ix77lNwS2cJN_lzQKPKLqLrsaAkrIbvSELvLyJU9B5mfv1IBClFA5cLAUc3v88LjuVqWKmKZjzP_SV1NQiHZT8n-iJQGX9HF"
$ws.Range("J6").Value = 4
$ws.Range("K6").Value = 0

# --- Page setup (print area / paper) ------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection, matching the authored state -----------------------------
$ws.Range("G6").Select()
